$d = $word.ActiveDocument
$t = $d.Tables.Item(2)

function Fill-Cell($cell, $text, $withBookmark) {
    $bmXml = ""
    if ($withBookmark) {
        $bmXml = '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'
    }
    $ooxml = @"
<?xml version="1.0" encoding="utf-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p w:rsidR="00780210" w:rsidRPr="00F93EB7" w:rsidRDefault="00780210" w:rsidP="00F93EB7">
<w:pPr>
<w:jc w:val="lowKashida"/>
<w:rPr>
<w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/>
<w:color w:val="000000" w:themeColor="text1"/>
<w:sz w:val="18"/>
<w:szCs w:val="18"/>
</w:rPr>
</w:pPr>
<w:r>
<w:rPr>
<w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/>
<w:color w:val="000000" w:themeColor="text1"/>
<w:sz w:val="18"/>
<w:szCs w:val="18"/>
</w:rPr>
<w:t>$text</w:t>
</w:r>
$bmXml
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
    $cell.Range.InsertXML($ooxml)
    $cell.Range.Paragraphs.Item(1).Range.Delete()
}

# Remove the _GoBack bookmark from the previous row's last cell (row 11, "Bagged Decision Trees")
$prevLastCell = $t.Cell(11, 12)
$ooxmlNoBm = @"
<?xml version="1.0" encoding="utf-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p w:rsidR="00780210" w:rsidRPr="00F93EB7" w:rsidRDefault="00AB3538" w:rsidP="00AB3538">
<w:pPr>
<w:jc w:val="lowKashida"/>
<w:rPr>
<w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/>
<w:color w:val="000000" w:themeColor="text1"/>
<w:sz w:val="18"/>
<w:szCs w:val="18"/>
</w:rPr>
</w:pPr>
<w:r w:rsidRPr="00AB3538">
<w:rPr>
<w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/>
<w:color w:val="000000" w:themeColor="text1"/>
<w:sz w:val="18"/>
<w:szCs w:val="18"/>
</w:rPr>
<w:t>0.999051</w:t>
</w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
$prevLastCell.Range.InsertXML($ooxmlNoBm)
$prevLastCell.Range.Paragraphs.Item(1).Range.Delete()

# Fill the "Bagged SVM" row (row 12) result cells
Fill-Cell $t.Cell(12, 3) "0.001666" $false
Fill-Cell $t.Cell(12, 4) "0.031755" $false
Fill-Cell $t.Cell(12, 5) "0.040816" $false
Fill-Cell $t.Cell(12, 6) "0.65152" $false
Fill-Cell $t.Cell(12, 7) "0.115697" $false
Fill-Cell $t.Cell(12, 8) "0.021914" $false
Fill-Cell $t.Cell(12, 9) "0.103569" $false
Fill-Cell $t.Cell(12, 10) "0.000986" $false
Fill-Cell $t.Cell(12, 11) "0.629221" $false
Fill-Cell $t.Cell(12, 12) "0.654758" $true
